# Updating changes related to env setup
# Bump the LiveSLR build number shown in the "Version" column (B2).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B2").Value = "Copyright @ 2023 Cytel Inc. LiveSLR 4.0.0.0 - Build #54694"
